# Auto-generated edit script applying the diff to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I header ("Other found locations")
$ws.Cells.Item(1, 9).Value = "Other found locations"

# Row 2
$ws.Cells.Item(2, 5).Value = "[ D.%D’Ardes%null%2,     F.%Santilli%null%1,     M. T.%Guagnano%null%1,     M.%Bucci%null%1,     F.%Cipollone%null%1,  D.%D’Ardes%null%0,  F.%Santilli%null%1,  M. T.%Guagnano%null%1,  M.%Bucci%null%1,  F.%Cipollone%null%1]"
$ws.Cells.Item(2, 9).Value = ""

# Row 3
$ws.Cells.Item(3, 5).Value = "[Chen%Chen%NULL%0,    Feng%Wang%NULL%0,    Peng%Chen%NULL%1,    Jiangang%Jiang%NULL%0,    Guanglin%Cui%NULL%1,    Ning%Zhou%NULL%0,    Francesco%Moroni%NULL%1,    Javid J.%Moslehi%NULL%0,    Enrico%Ammirati%enrico.ammirati@ospedaleniguarda.it%0,    Dao Wen%Wang%dwwang@tjh.tjmu.edu.cn%0]"
$ws.Cells.Item(3, 9).Value = "_PMC"

# Row 4
$ws.Cells.Item(4, 5).Value = "[Marcello%Covino%NULL%0,    Giuseppe%De Matteis%giuseppe.dematteis@policlinicogemelli.it%1,    Maria Livia%Burzo%NULL%2,    Maria Livia%Burzo%NULL%0,    Michele%Santoro%NULL%1,    Mariella%Fuorlo%NULL%1,    Luca%Sabia%NULL%1,    Claudio%Sandroni%NULL%1,    Antonio%Gasbarrini%NULL%0,    Francesco%Franceschi%NULL%1,    Giovanni%Gambassi%NULL%1,    NULL%NULL%NULL%0]"
$ws.Cells.Item(4, 9).Value = "_PMC"

# Row 5
$ws.Cells.Item(5, 5).Value = "[C%Felice%NULL%1,    C%Nardin%NULL%1,    G L%Di Tanna%NULL%1,    U%Grossi%NULL%1,    E%Bernardi%NULL%1,    L%Scaldaferri%NULL%1,    M%Romagnoli%NULL%1,    L%Tonon%NULL%1,    P%Cavasin%NULL%1,    S%Novello%NULL%1,    R%Scarpa%NULL%1,    A%Farnia%NULL%1,    E%De Menis%NULL%1,    R%Rigoli%NULL%1,    F%Cinetto%NULL%1,    P%Pauletto%NULL%1,    C%Agostini%NULL%1,    M%Rattazzi%marcello.rattazzi@unipd.it%1]"
$ws.Cells.Item(5, 9).Value = "_PMC"

# Row 6
$ws.Cells.Item(6, 5).Value = "[Zheyong%Huang%NULL%0,    Jiatian%Cao%NULL%1,    Yumeng%Yao%NULL%1,    Xuejuan%Jin%NULL%1,    Zhe%Luo%NULL%1,    Yuan%Xue%NULL%1,    Chouwen%Zhu%NULL%1,    Yanan%Song%NULL%1,    Ying%Wang%NULL%0,    Yunzeng%Zou%NULL%1,    Juying%Qian%NULL%1,    Kaihuan%Yu%NULL%1,    Hui%Gong%NULL%1,    Junbo%Ge%NULL%1]"
$ws.Cells.Item(6, 9).Value = "_PMC"

# Row 7
$ws.Cells.Item(7, 5).Value = "[Katherine W%Lam%NULL%0,    Kenneth W%Chow%NULL%1,    Jonathan%Vo%NULL%1,    Wei%Hou%NULL%0,    Haifang%Li%NULL%0,    Paul S%Richman%NULL%1,    Sandeep K%Mallipattu%NULL%1,    Hal A%Skopicki%NULL%1,    Adam J%Singer%NULL%1,    Tim Q%Duong%tim.duong@stonybrook.edu%1]"
$ws.Cells.Item(7, 9).Value = "_PMC"

# Row 8
$ws.Cells.Item(8, 5).Value = "[Xiulan%Liu%NULL%0,    Yi%Liu%NULL%0,    Keliang%Chen%NULL%1,    Suying%Yan%NULL%1,    Xiangrong%Bai%NULL%1,    Juan%Li%947281063@qq.com%0,    Dong%Liu%ld2069@outlook.com%0,    Dong%Liu%ld2069@outlook.com%0]"
$ws.Cells.Item(8, 9).Value = "_PMC"

# Row 9
$ws.Cells.Item(9, 5).Value = "[ Yasushi%Matsuzawa%null%0,     Hisao%Ogawa%null%1,     Kazuo%Kimura%null%1,     Masaaki%Konishi%null%1,     Jin%Kirigaya%null%1,     Kazuki%Fukui%null%1,     Kengo%Tsukahara%null%1,     Hiroyuki%Shimizu%null%0,     Keisuke%Iwabuchi%null%0,     Yu%Yamada%null%1,     Kenichiro%Saka%null%1,     Ichiro%Takeuchi%null%0,     Toshio%Hirano%null%1,     Kouichi%Tamura%null%1,  Yasushi%Matsuzawa%null%0,  Hisao%Ogawa%null%1,  Kazuo%Kimura%null%1,  Masaaki%Konishi%null%1,  Jin%Kirigaya%null%1,  Kazuki%Fukui%null%1,  Kengo%Tsukahara%null%1,  Hiroyuki%Shimizu%null%1,  Keisuke%Iwabuchi%null%1,  Yu%Yamada%null%1,  Kenichiro%Saka%null%1,  Ichiro%Takeuchi%null%1,  Toshio%Hirano%null%1,  Kouichi%Tamura%null%1]"
$ws.Cells.Item(9, 9).Value = ""

# Row 10
$ws.Cells.Item(10, 5).Value = "[Juan%Meng%NULL%0,    Guohui%Xiao%NULL%1,    Juanjuan%Zhang%NULL%1,    Xing%He%NULL%1,    Min%Ou%NULL%1,    Jing%Bi%NULL%1,    Rongqing%Yang%NULL%1,    Wencheng%Di%NULL%1,    Zhaoqin%Wang%NULL%0,    Zigang%Li%NULL%1,    Hong%Gao%NULL%0,    Lei%Liu%NULL%0,    Guoliang%Zhang%NULL%1]"
$ws.Cells.Item(10, 9).Value = "_PMC"

# Row 11
$ws.Cells.Item(11, 5).Value = "[ Lili%Niu%null%2,     Yanchen%Guo%null%1,     Zhengrong%Lin%null%1,     Zhe%Shi%null%1,     Tianyuan%Bian%null%1,     Lin%Qi%null%1,     Long%Meng%null%1,     Anthony A.%Grace%null%1,     Hairong%Zheng%null%1,     Ti-Fei%Yuan%null%1,  Lili%Niu%null%0,  Yanchen%Guo%null%1,  Zhengrong%Lin%null%1,  Zhe%Shi%null%1,  Tianyuan%Bian%null%1,  Lin%Qi%null%1,  Long%Meng%null%1,  Anthony A.%Grace%null%1,  Hairong%Zheng%null%1,  Ti-Fei%Yuan%null%1]"
$ws.Cells.Item(11, 9).Value = ""

# Row 12
$ws.Cells.Item(12, 9).Value = ""

# Row 13
$ws.Cells.Item(13, 5).Value = "[Celestino%Sardu%drsarducele@gmail.com%0,    Paolo%Maggi%NULL%2,    Paolo%Maggi%NULL%0,    Vincenzo%Messina%NULL%1,    Pasquale%Iuliano%NULL%2,    Pasquale%Iuliano%NULL%0,    Antonio%Sardu%NULL%1,    Vincenzo%Iovinella%NULL%1,    Giuseppe%Paolisso%NULL%1,    Raffaele%Marfella%NULL%1]"
$ws.Cells.Item(13, 9).Value = "_PMC"

# Row 14
$ws.Cells.Item(14, 5).Value = "[International%Society of Hypertension%coreGivesNoEmail%1,   Li%J%coreGivesNoEmail%1,   Zhang%P%coreGivesNoEmail%1]"
$ws.Cells.Item(14, 6).Value = "not found"
$ws.Cells.Item(14, 7).Value = "N/A"
$ws.Cells.Item(14, 9).Value = ""

# Row 15
$ws.Cells.Item(15, 5).Value = "[Naci%Şenkal%NULL%0,    Rasimcan%Meral%NULL%1,    Alpay%Medetalibeyoğlu%NULL%1,    Hilal%Konyaoğlu%NULL%1,    Murat%Köse%NULL%1,    Tufan%Tükek%NULL%1]"
$ws.Cells.Item(15, 9).Value = "_PMC"

# Row 16
$ws.Cells.Item(16, 9).Value = ""

# Row 17
$ws.Cells.Item(17, 5).Value = "[Zhongchao%Wang%NULL%1,    Dewei%Zhang%NULL%1,    Shengming%Wang%NULL%1,    Yanhua%Jin%NULL%1,    Jianbo%Huan%NULL%1,    Yue%Wu%NULL%0,    Cheng%Xia%NULL%1,    Zhe%Li%NULL%0,    Xingshun%Qi%NULL%1,    Duanzhen%Zhang%NULL%1,    Xiumin%Han%NULL%1,    Xianyang%Zhu%NULL%1,    Ying%Qu%NULL%1,    Qiguang%Wang%NULL%1]"
$ws.Cells.Item(17, 9).Value = "_PMC"

# Row 18
$ws.Cells.Item(18, 5).Value = "[Honglu%Zhou%NULL%1,    Songmei%Wang%NULL%1,    Lorenz%von Seidlein%NULL%1,    Xuanyi%Wang%xywang@shmu.edu.cn%1]"
$ws.Cells.Item(18, 9).Value = "_PMC_Springer"

# Row 19
$ws.Cells.Item(19, 9).Value = ""

# Row 20
$ws.Cells.Item(20, 5).Value = "[Yuan%Yuan%NULL%0,    Dan%Liu%NULL%0,    Shaoqing%Zeng%NULL%1,    Siyuan%Wang%NULL%1,    Sen%Xu%NULL%0,    Ya%Wang%NULL%0,    Ruidi%Yu%NULL%1,    Yue%Gao%NULL%0,    Huayi%Li%NULL%1,    Xinxia%Feng%NULL%1,    Ning%Zhou%NULL%0,    Chunxia%Zhao%NULL%1,    Qinglei%Gao%NULL%1]"
$ws.Cells.Item(20, 9).Value = "_PMC_elsevier"

# Row 21
$ws.Cells.Item(21, 5).Value = "[Peng%Zhang%NULL%0,    Lihua%Zhu%NULL%0,    Jingjing%Cai%NULL%0,    Fang%Lei%NULL%0,    Juan-Juan%Qin%NULL%0,    Jing%Xie%NULL%2,    Ye-Mao%Liu%NULL%0,    Yan-Ci%Zhao%NULL%1,    Xuewei%Huang%NULL%0,    Lijin%Lin%NULL%0,    Meng%Xia%NULL%0,    Ming-Ming%Chen%NULL%0,    Xu%Cheng%NULL%1,    Xiao%Zhang%NULL%0,    Deliang%Guo%NULL%1,    Yuanyuan%Peng%NULL%1,    Yan-Xiao%Ji%NULL%0,    Jing%Chen%NULL%0,    Zhi-Gang%She%NULL%0,    Yibin%Wang%NULL%0,    Qingbo%Xu%NULL%0,    Renfu%Tan%NULL%1,    Haitao%Wang%NULL%1,    Jun%Lin%NULL%2,    Pengcheng%Luo%NULL%1,    Shouzhi%Fu%NULL%0,    Hongbin%Cai%NULL%1,    Ping%Ye%NULL%0,    Bing%Xiao%NULL%1,    Weiming%Mao%NULL%0,    Liming%Liu%NULL%1,    Youqin%Yan%NULL%2,    Mingyu%Liu%NULL%1,    Manhua%Chen%NULL%1,    Xiao-Jing%Zhang%NULL%0,    Xinghuan%Wang%NULL%0,    Rhian M.%Touyz%NULL%1,    Jiahong%Xia%NULL%0,    Bing-Hong%Zhang%NULL%0,    Xiaodong%Huang%NULL%0,    Yufeng%Yuan%NULL%0,    Rohit%Loomba%NULL%1,    Peter P.%Liu%NULL%1,    Hongliang%Li%NULL%0]"
$ws.Cells.Item(21, 9).Value = "_PMC"

# Row 22
$ws.Cells.Item(22, 5).Value = "[Xian%Zhou%NULL%1,    Jingkang%Zhu%NULL%1,    Tao%Xu%NULL%0]"
$ws.Cells.Item(22, 9).Value = "_PMC"

